$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Brn_Code (column C) was stored as the text "195"; switch it to the
# numeric value 19 for rows 2-4.
$ws.Range("C2:C4").Value = 19

# Classification_TYPE (column J) was blank; set the processed flags.
$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 0
$ws.Range("J4").Value = 0

# The author's last selection moved from J22 to J14.
$ws.Range("J14").Select()
